$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: ", прави проследяването на проблемие, които са били " ->
#           ", прави проследяването на проблемите, които са били "
#           split into 3 runs: "...проблеми" | "т" | "е, които са били "
#
# Splitting a run into separate <w:r> elements (while keeping identical
# rPr) is achieved by briefly adding a bookmark at the desired boundary
# and then deleting it again -- Word breaks the run there and the split
# survives the bookmark's removal.
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("проблемие, които са били", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start1 = $r1.Start
# Insert the missing "т" right before the final "е" of "проблемие" (index 8).
$insPoint = $d.Range($start1 + 8, $start1 + 8)
$insPoint.InsertAfter("т")
# Now the text reads "...проблемите, които са били ..." in a single run.
# Split it into three runs at the boundaries around the inserted "т".
$sp1 = $d.Range($start1 + 8, $start1 + 8)
$d.Bookmarks.Add("TmpSplitA1", $sp1) | Out-Null
$d.Bookmarks("TmpSplitA1").Delete()

$sp2 = $d.Range($start1 + 9, $start1 + 9)
$d.Bookmarks.Add("TmpSplitA2", $sp2) | Out-Null
$d.Bookmarks("TmpSplitA2").Delete()

# ---------------------------------------------------------------------
# Change 2: rewrite the "Closure phase" conclusion sentence.
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("-а ще може да намира по-лесно евентуални бъгове, тъй като той ще се е поучил от вече намерените", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start2 = $r2.Start
$r2.Text = "-ите, които са прегледали документите със заключения, ще може да намират по-лесно евентуални бъгове, тъй като ще има написани поуки от вече намерени бъгове"

# Split the merged run into the 7 runs seen in the target, at the
# boundaries between each textual chunk.
$spB1 = $d.Range($start2 + 53, $start2 + 53)
$d.Bookmarks.Add("TmpSplitB1", $spB1) | Out-Null
$d.Bookmarks("TmpSplitB1").Delete()

$spB2 = $d.Range($start2 + 71, $start2 + 71)
$d.Bookmarks.Add("TmpSplitB2", $spB2) | Out-Null
$d.Bookmarks("TmpSplitB2").Delete()

$spB3 = $d.Range($start2 + 72, $start2 + 72)
$d.Bookmarks.Add("TmpSplitB3", $spB3) | Out-Null
$d.Bookmarks("TmpSplitB3").Delete()

$spB4 = $d.Range($start2 + 99, $start2 + 99)
$d.Bookmarks.Add("TmpSplitB4", $spB4) | Out-Null
$d.Bookmarks("TmpSplitB4").Delete()

$spB5 = $d.Range($start2 + 110, $start2 + 110)
$d.Bookmarks.Add("TmpSplitB5", $spB5) | Out-Null
$d.Bookmarks("TmpSplitB5").Delete()

# The last boundary (offset 148) becomes the relocated "_GoBack" bookmark;
# adding a bookmark with that reserved name automatically removes/moves the
# one that previously sat near the end of the document.
$goBackPoint = $d.Range($start2 + 148, $start2 + 148)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

# ---------------------------------------------------------------------
# Change 3: "if (false) {" -> "if" | " (false) {" (two runs)
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("if (false) {", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start3 = $r3.Start
$spC1 = $d.Range($start3 + 2, $start3 + 2)
$d.Bookmarks.Add("TmpSplitC1", $spC1) | Out-Null
$d.Bookmarks("TmpSplitC1").Delete()
